# banco_girias.xlsx - add new slang entries (rows 3-11) to Sheet 1
# Columns: A=id_giria B=titulo C=traducaoTitulo D=descricao E=exemplo F=exemploTraducao G=dificuldade

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Tosser
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Tosser"
$ws.Range("C3").Value = "Sem tradução"
$ws.Range("D3").Value = "Idiota"
$ws.Range("E3").Value = "####"
$ws.Range("F3").Value = "#####"
$ws.Range("G3").Value = 1

# Row 4 - cock-up
$ws.Range("B4").Value = "cock-up"
$ws.Range("C4").Value = "Sem tradução"
$ws.Range("D4").Value = "arruinar, estragar algo"
$ws.Range("E4").Value = "we've made a total cock-up of it"
$ws.Range("F4").Value = "nós fizemos uma confusão total"

# Row 5 - Bloody
$ws.Range("B5").Value = "Bloody"
$ws.Range("C5").Value = "sangrento"
$ws.Range("D5").Value = "Para algo que não deu certo"
$ws.Range("E5").Value = "Bloody!"
$ws.Range("F5").Value = "Que Droga!, Que saco!"

# Row 6 - Give You A Bell
$ws.Range("B6").Value = "Give You A Bell"
$ws.Range("C6").Value = "Dê-lhe um sino"
$ws.Range("D6").Value = "ligar para voce"
$ws.Range("E6").Value = "####"
$ws.Range("F6").Value = "####"

# Row 7 - Blimey!
$ws.Range("B7").Value = "Blimey!"
$ws.Range("C7").Value = "Caramba!"
$ws.Range("D7").Value = "Usado para expressar surpresa"
$ws.Range("E7").Value = "Blimey! did you see that mole rat lookin brit's teeth?"
$ws.Range("F7").Value = "Caramba! Você viu esse rato procurando os dentes do britânico?"

# Row 8 - Wanker
$ws.Range("B8").Value = "Wanker"
$ws.Range("C8").Value = "Sem tradução"
$ws.Range("D8").Value = "Idiota"
$ws.Range("E8").Value = "#####"
$ws.Range("F8").Value = "####"

# Row 9 - Gutted
$ws.Range("B9").Value = "Gutted"
$ws.Range("C9").Value = "Esvaziado"
$ws.Range("D9").Value = "Devastado(a)"
$ws.Range("E9").Value = "After getting a few fish each, they swam in the pond before they went back to the beach to clean and gut the fish and prepare them for dinner."
$ws.Range("F9").Value = "Depois de pegar alguns peixes, eles nadaram na lagoa antes de voltarem para a praia para limpar e destruir o peixe e prepará-los para o jantar."

# Row 10 - Bespoke
$ws.Range("B10").Value = "Bespoke "
$ws.Range("C10").Value = "Sem tradução"
$ws.Range("D10").Value = "Feito sob medida"
$ws.Range("E10").Value = "a bespoke suit"
$ws.Range("F10").Value = "um terno feito sob medida"

# Row 11 - Chuffed
$ws.Range("B11").Value = "Chuffed"
$ws.Range("C11").Value = "Sem tradução"
$ws.Range("D11").Value = "orgulhoso(a)"
$ws.Range("E11").Value = "I'm dead chuffed to have wo"
$ws.Range("F11").Value = "Estou morta e feliz por ter ganhado"

# Update view / selection and print setup to match the edited workbook
$ws.Range("C13").Select()
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
